$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.099699974060059
$ws.Range("B1").Value = 1.275647401809692
$ws.Range("C1").Value = 1.618452787399292
$ws.Range("D1").Value = 3.185163736343384
$ws.Range("E1").Value = 15
